# SDVs.xlsx working-model consolidation edit:
# Rename the "sprocketAngVel" parameter to "sprocketRPM" throughout the
# workbook (the cell that holds the label on Sheet1, the workbook-level
# defined name that still points at Sheet1!$C$6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B6 holds the text label for the row - update it in place.
$ws.Range("B6").Value = "sprocketRPM"

# The workbook defined name keeps referring to the same cell (Sheet1!$C$6),
# only its name changes.
$wb.Names.Item("sprocketAngVel").Name = "sprocketRPM"
